# nexial-macro.xlsx : add two new "web" commands to the #system lookup list
#   - assertSelectOptionsAbsent(locator,options)
#   - assertSelectOptionsPresent(locator,options)
#
# The hidden "#system" sheet stores, per-command-category, an alphabetically
# sorted list of command signatures in its own column; the named range
# "web" (column AE) backs the MacroLibrary sheet's drop-down validation for
# the "web" category. The two new commands sort alphabetically right before
# the existing "assertSingleSelect(locator)" entry (currently row 39), so
# every entry from row 39 down shifts two rows later and the named range's
# extent grows from $AE$2:$AE$156 to $AE$2:$AE$158.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Shift AE39:AE156 down to AE41:AE158, leaving two fresh blank cells at
# AE39:AE40. Using a column-only range (not EntireRow) so the other
# category columns on this sheet (F, I, T, AD, ...) are left untouched.
$ws.Range("AE39:AE40").Insert()

# Populate the two freshly inserted cells with the new command signatures,
# in alphabetical order.
$ws.Range("AE39").Value = "assertSelectOptionsAbsent(locator,options)"
$ws.Range("AE40").Value = "assertSelectOptionsPresent(locator,options)"

# Extend the "web" defined name to cover the two extra rows.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$AE`$2:`$AE`$158"
    }
}
